# Apply the thesis-notes edit:
#  - add w:proofErr spell/grammar markers around a few existing runs
#    (splitting single runs into multiple where the marker needs to sit
#    mid-text), and
#  - append two new paragraphs ("Calibration-Free Localization.m >>" and
#    the "set rand factor for gamma ..." note) before the closing
#    bookmark paragraph.
#
# w:proofErr markers aren't reachable through any higher-level Range/Find
# property, so we drop down to Range.InsertXML, which *replaces* the
# content of the (non-collapsed) range it's called on with the supplied
# WordprocessingML. We work on one paragraph Range at a time, and do the
# later paragraphs first so earlier paragraphs' indices/offsets stay
# valid while we still need them.

$d = $word.ActiveDocument
$pkgOpen = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body>"
$pkgClose = "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

# --- Paragraphs 6 & 7 ("Root mean square ... EKF" + the trailing empty
#     paragraph that held the _GoBack bookmark) -> becomes: the EKF
#     paragraph (bookmark removed), a blank paragraph, the new
#     "Calibration-Free Localization.m >>" paragraph, the new "set rand
#     factor for gamma" paragraph, a blank paragraph, and finally the
#     bookmark-only paragraph.
$p6 = $d.Paragraphs(6)
$p7 = $d.Paragraphs(7)
$combined = $d.Range($p6.Range.Start, $p7.Range.End)
$body = @"
<w:p>
<w:r><w:t>Root mean square</w:t></w:r>
<w:r><w:t xml:space="preserve"> of derivation from the true trajectory of EKF</w:t></w:r>
</w:p>
<w:p/>
<w:p>
<w:r><w:t xml:space="preserve">Calibration-Free </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>Localization</w:t></w:r>
<w:r><w:t>.m</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> &gt;&gt;</w:t></w:r>
</w:p>
<w:p>
<w:pPr>
<w:autoSpaceDE w:val="0"/>
<w:autoSpaceDN w:val="0"/>
<w:adjustRightInd w:val="0"/>
<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
<w:ind w:firstLine="720"/>
<w:rPr>
<w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
<w:sz w:val="24"/>
<w:szCs w:val="24"/>
</w:rPr>
</w:pPr>
<w:proofErr w:type="gramStart"/>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
<w:color w:val="228B22"/>
<w:sz w:val="26"/>
<w:szCs w:val="26"/>
</w:rPr>
<w:t>set</w:t>
</w:r>
<w:proofErr w:type="gramEnd"/>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
<w:color w:val="228B22"/>
<w:sz w:val="26"/>
<w:szCs w:val="26"/>
</w:rPr>
<w:t xml:space="preserve"> rand factor for gamma</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
<w:color w:val="228B22"/>
<w:sz w:val="26"/>
<w:szCs w:val="26"/>
</w:rPr>
<w:t xml:space="preserve"> seems could solve local minimal problem :D</w:t>
</w:r>
</w:p>
<w:p/>
<w:p>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
"@
$combined.InsertXML($pkgOpen + $body + $pkgClose) | Out-Null

# --- Paragraph 4 ("Plot of errors(mean variance) to #nodes actived in
#     the last n seconds") -> split into runs so gramStart/gramEnd can
#     bracket "errors(" and spellStart/spellEnd can bracket "actived".
$p4 = $d.Paragraphs(4)
$body = @"
<w:p>
<w:r><w:t xml:space="preserve">Plot of </w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:t>errors(</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:t xml:space="preserve">mean variance) to #nodes </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>actived</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> in the last n seconds </w:t></w:r>
</w:p>
"@
$p4.Range.InsertXML($pkgOpen + $body + $pkgClose) | Out-Null

# --- Paragraph 2 ("Mad()") -> wrap with gramStart/gramEnd.
$p2 = $d.Paragraphs(2)
$body = @"
<w:p>
<w:proofErr w:type="gramStart"/>
<w:r><w:t>Mad</w:t></w:r>
<w:r><w:t>()</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
</w:p>
"@
$p2.Range.InsertXML($pkgOpen + $body + $pkgClose) | Out-Null

# --- Paragraph 1 ("Trimmean()") -> wrap with spellStart/spellEnd around
#     "Trimmean" and gramStart/gramEnd around the whole "Trimmean()".
$p1 = $d.Paragraphs(1)
$body = @"
<w:p>
<w:proofErr w:type="spellStart"/>
<w:proofErr w:type="gramStart"/>
<w:r><w:t>Trimmean</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t>()</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
</w:p>
"@
$p1.Range.InsertXML($pkgOpen + $body + $pkgClose) | Out-Null
